$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 was previously stored as an inline string "456"; it should be a real number.
$ws.Range("C2").Value = 456

# Append new data row 3: 190, 119, 456, 2023-05-04 (serial 45050), 0, "190-119-text"
$ws.Range("A3").Value = 190
$ws.Range("B3").Value = 119
$ws.Range("C3").Value = 456
$ws.Range("D3").Value = 45050
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "190-119-text"
